$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows above the current row 350 (Early Majestic / Especial).
# This pushes the existing "Early Majestic" rows (350-352) down to (354-356)
# unchanged, and leaves 4 fresh blank rows at 350-353 for the new weekly data.
$ws.Rows.Item(350).Insert()
$ws.Rows.Item(350).Insert()
$ws.Rows.Item(350).Insert()
$ws.Rows.Item(350).Insert()

# Common values shared by all rows in this block.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$semana = 4
$rubro = "Fruta"
$grupoId = 100103
$grupo = "Frutos de hueso (carozo)"
$especieId = 100103004
$especie = "Durazno"
$unidad = "`$/bins (400 kilos)"
$kilos = 400

function Set-Row($r, $fecha, $variedad, $calidad, $calibre, $precioMin, $precioMax, $precioProm, $provincia, $precioKilo) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $semana
    $ws.Cells.Item($r, 6).Value = $rubro
    $ws.Cells.Item($r, 7).Value = $grupoId
    $ws.Cells.Item($r, 8).Value = $grupo
    $ws.Cells.Item($r, 9).Value = $especieId
    $ws.Cells.Item($r, 10).Value = $especie
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $calibre
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $provincia
    $ws.Cells.Item($r, 19).Value = $precioKilo
    $ws.Cells.Item($r, 20).Value = $kilos
}

# New week (2022-02-03, serial 44595) data.
Set-Row 350 44595 "Carson"       "Primera" 20 315000 320000 317500 "Región de O'Higgins" 794
Set-Row 351 44595 "Carson"       "Segunda" 20 275000 280000 277500 "Región de O'Higgins" 694
Set-Row 352 44595 "Doctor Davis" "Primera" 20 315000 320000 317500 "Región de O'Higgins" 794
Set-Row 353 44595 "Doctor Davis" "Segunda" 20 285000 290000 287500 "Región de O'Higgins" 719

"Done"
